# Auto-generated Excel COM-interop script to apply crypto list update
# (cryptos.xlsx price/volume refresh, commit: Fri Dec  8 07:35:25 UTC 2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '43.261.37'
$ws.Cells.Item(2, 4).Style = 'Normal'
$ws.Cells.Item(2, 5).Value = '  -1.55%  '

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = '@'
$ws.Cells.Item(3, 4).Value = '2.362.73'
$ws.Cells.Item(3, 4).Style = 'Normal'
$ws.Cells.Item(3, 5).Value = '  +4.24%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.12%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '233.29'
$ws.Cells.Item(5, 4).Style = 'Normal'
$ws.Cells.Item(5, 5).Value = '  -0.19%  '

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '0.647'
$ws.Cells.Item(6, 4).Style = 'Normal'
$ws.Cells.Item(6, 5).Value = '  -1.45%  '

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '71.97'
$ws.Cells.Item(7, 4).Style = 'Normal'
$ws.Cells.Item(7, 5).Value = '  +12.81%  '

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '0.476'
$ws.Cells.Item(9, 4).Style = 'Normal'
$ws.Cells.Item(9, 5).Value = '  +5.66%  '

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '0.0969'
$ws.Cells.Item(10, 4).Style = 'Normal'
$ws.Cells.Item(10, 5).Value = '  -0.82%  '

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '56.84'
$ws.Cells.Item(11, 4).Style = 'Normal'
$ws.Cells.Item(11, 5).Value = '  -1.82%  '

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '27.09'
$ws.Cells.Item(12, 4).Style = 'Normal'
$ws.Cells.Item(12, 5).Value = '  +1.87%  '

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '2.724.29'
$ws.Cells.Item(13, 4).Style = 'Normal'
$ws.Cells.Item(13, 5).Value = '  +4.76%  '

# Row 14
$ws.Cells.Item(14, 5).Value = '  +0.61%  '

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '15.91'
$ws.Cells.Item(15, 4).Style = 'Normal'
$ws.Cells.Item(15, 5).Value = '  +1.72%  '

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '6.26'
$ws.Cells.Item(16, 4).Style = 'Normal'
$ws.Cells.Item(16, 5).Value = '  +1.96%  '

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '0.855'
$ws.Cells.Item(17, 4).Style = 'Normal'
$ws.Cells.Item(17, 5).Value = '  +1.65%  '

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '2.366.90'
$ws.Cells.Item(18, 4).Style = 'Normal'
$ws.Cells.Item(18, 5).Value = '  +4.30%  '

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '43.220.89'
$ws.Cells.Item(19, 4).Style = 'Normal'
$ws.Cells.Item(19, 5).Value = '  -1.42%  '

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '0.0₃0989'
$ws.Cells.Item(20, 4).Style = 'Normal'
$ws.Cells.Item(20, 5).Value = '  +0.71%  '

# Row 21
$ws.Cells.Item(21, 2).Value = 'Uniswap'
$ws.Cells.Item(21, 3).Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '6.31'
$ws.Cells.Item(21, 4).Style = 'Normal'
$ws.Cells.Item(21, 5).Value = '  +2.32%  '

# Row 22
$ws.Cells.Item(22, 2).Value = 'Litecoin'
$ws.Cells.Item(22, 3).Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '74.19'
$ws.Cells.Item(22, 4).Style = 'Normal'
$ws.Cells.Item(22, 5).Value = '  +0.51%  '

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '248.82'
$ws.Cells.Item(23, 4).Style = 'Normal'
$ws.Cells.Item(23, 5).Value = '  -0.43%  '

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '0.999'
$ws.Cells.Item(24, 4).Style = 'Normal'
$ws.Cells.Item(24, 5).Value = '  -0.04%  '

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '3.69'
$ws.Cells.Item(25, 4).Style = 'Normal'
$ws.Cells.Item(25, 5).Value = '  +1.21%  '

# Row 26
$ws.Cells.Item(26, 5).Value = '  -0.16%  '

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '2.26'
$ws.Cells.Item(27, 4).Style = 'Normal'
$ws.Cells.Item(27, 5).Value = '  +0.91%  '

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '9.93'
$ws.Cells.Item(28, 4).Style = 'Normal'
$ws.Cells.Item(28, 5).Value = '  +0.40%  '

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '22.28'
$ws.Cells.Item(29, 4).Style = 'Normal'
$ws.Cells.Item(29, 5).Value = '  +1.32%  '

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '173.70'
$ws.Cells.Item(30, 4).Style = 'Normal'
$ws.Cells.Item(30, 5).Value = '  -0.22%  '

# Row 31
$ws.Cells.Item(31, 5).Value = '  +5.36%  '

# Row 32
$ws.Cells.Item(32, 5).Value = '  -5.65%  '

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '0.126'
$ws.Cells.Item(33, 4).Style = 'Normal'
$ws.Cells.Item(33, 5).Value = '  -0.71%  '

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '4.95'
$ws.Cells.Item(34, 4).Style = 'Normal'
$ws.Cells.Item(34, 5).Value = '  -1.00%  '

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '0.0689'
$ws.Cells.Item(35, 4).Style = 'Normal'
$ws.Cells.Item(35, 5).Value = '  +0.54%  '

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '5.02'
$ws.Cells.Item(36, 4).Style = 'Normal'
$ws.Cells.Item(36, 5).Value = '  +1.06%  '

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '6.56'
$ws.Cells.Item(37, 4).Style = 'Normal'
$ws.Cells.Item(37, 5).Value = '  +1.95%  '

# Row 38
$ws.Cells.Item(38, 2).Value = 'LidoDAOToken'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '2.44'
$ws.Cells.Item(38, 4).Style = 'Normal'
$ws.Cells.Item(38, 5).Value = '  +5.91%  '

# Row 39
$ws.Cells.Item(39, 2).Value = 'RenderToken'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '3.68'
$ws.Cells.Item(39, 4).Style = 'Normal'
$ws.Cells.Item(39, 5).Value = '  -0.77%  '

# Row 40
$ws.Cells.Item(40, 5).Value = '  -0.73%  '

# Row 41
$ws.Cells.Item(41, 5).Value = '  +0.12%  '

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '8.87'
$ws.Cells.Item(42, 4).Style = 'Normal'
$ws.Cells.Item(42, 5).Value = '  +1.16%  '

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '18.40'
$ws.Cells.Item(43, 4).Style = 'Normal'
$ws.Cells.Item(43, 5).Value = '  +5.99%  '

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '1.17'
$ws.Cells.Item(44, 4).Style = 'Normal'
$ws.Cells.Item(44, 5).Value = '  +8.20%  '

# Row 45
$ws.Cells.Item(45, 2).Value = 'FTXToken'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '4.52'
$ws.Cells.Item(45, 4).Style = 'Normal'
$ws.Cells.Item(45, 5).Value = '  -1.07%  '

# Row 46
$ws.Cells.Item(46, 2).Value = 'Aave'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '99.48'
$ws.Cells.Item(46, 4).Style = 'Normal'
$ws.Cells.Item(46, 5).Value = '  +0.75%  '

# Row 47
$ws.Cells.Item(47, 5).Value = '  +1.15%  '

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '0.0950'
$ws.Cells.Item(48, 4).Style = 'Normal'
$ws.Cells.Item(48, 5).Value = '  -0.32%  '

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '1.441.36'
$ws.Cells.Item(49, 4).Style = 'Normal'
$ws.Cells.Item(49, 5).Value = '  -1.15%  '

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '2.590.74'
$ws.Cells.Item(50, 4).Style = 'Normal'
$ws.Cells.Item(50, 5).Value = '  +4.66%  '

# Row 51
$ws.Cells.Item(51, 5).Value = '  -2.76%  '
